$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric Class values (-1 / 1) in column A (rows 2-16)
# with their text labels: -1 -> "Rejected", 1 -> "Accepted".
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq -1) {
        $cell.Value = "Rejected"
    } elseif ($cell.Value2 -eq 1) {
        $cell.Value = "Accepted"
    }
}

# Update the active selection to match the saved view state.
$ws.Range("B16").Select()
